$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column E values for rows 50..101 (small revisions) ---
$eUpdates = @{
    50 = 28608.8
    51 = 28626.1
    52 = 28386.3
    53 = 28115.5
    54 = 27879.2
    55 = 27710.3
    56 = 28020.6
    57 = 28376.5
    58 = 28471
    59 = 29442.5
    60 = 30132.3
    61 = 30509.6
    62 = 30891.4
    63 = 31321.3
    64 = 31435.1
    65 = 32044.5
    66 = 32597.1
    67 = 33088.9
    68 = 33294.3
    69 = 33533.5
    70 = 33923.7
    71 = 34447.6
    72 = 34741.9
    73 = 34815.5
    74 = 34867.3
    75 = 35031
    76 = 35002.3
    77 = 35405.8
    78 = 35499.5
    79 = 35964.7
    80 = 35830
    81 = 36180.5
    82 = 36417.8
    83 = 36341.3
    84 = 36473.4
    85 = 36568.9
    86 = 36267.9
    87 = 36646
    88 = 37303.2
    89 = 37635.1
    90 = 38054.8
    91 = 38386.3
    92 = 38272.3
    93 = 38706.1
    94 = 38598.5
    95 = 39019.2
    96 = 39330.9
    97 = 37817.8
    98 = 38551.6
    99 = 33665.2
    100 = 35492.8
    101 = 37787.9
}

foreach ($row in $eUpdates.Keys) {
    $ws.Range("E$row").Value = $eUpdates[$row]
}

# --- Update row 102 (01-01-2021) with revised figures ---
$ws.Range("B102").Value = 268019.7
$ws.Range("C102").Value = 56433.4
$ws.Range("D102").Value = 38075.6
$ws.Range("E102").Value = 39059.8
$ws.Range("F102").Value = 4087.3
$ws.Range("G102").Value = 3600.4
$ws.Range("H102").Value = 39367.9
$ws.Range("I102").Value = 29624.4
$ws.Range("J102").Value = 7786.9
$ws.Range("K102").Value = 20.7
$ws.Range("L102").Value = 10937.2
$ws.Range("M102").Value = 12265.4
$ws.Range("N102").Value = 125.3

# --- Add new row 103 (01-04-2021) ---
# Use a leading apostrophe so Excel stores this as text (matching the other
# "Serie" labels in column A) instead of auto-converting it to a date serial.
$ws.Range("A103").Value = "'01-04-2021"
# Re-apply the plain (unstyled) format used by the rest of column A so the
# quote-prefix formatting introduced above doesn't linger on the cell.
$ws.Range("A102").Copy()
$ws.Range("A103").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B103").Value = 292467.7
$ws.Range("C103").Value = 57426.7
$ws.Range("D103").Value = 39676.6
$ws.Range("E103").Value = 39460.4
$ws.Range("F103").Value = 4145.2
$ws.Range("G103").Value = 3775.5
$ws.Range("H103").Value = 41666.6
$ws.Range("I103").Value = 32833.5
$ws.Range("J103").Value = 8217.6
$ws.Range("K103").Value = 21
$ws.Range("L103").Value = 10700.1
$ws.Range("M103").Value = 12670.9
$ws.Range("N103").Value = 130.1

$wb.Save()
